$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.397.38"
$ws.Range("E2").Value = "  +0.64%  "

$ws.Range("D3").Value = "3.508.06"
$ws.Range("E3").Value = "  +3.84%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "668.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.50"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.431"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.11%  "

$ws.Range("E9").Value = "  +2.34%  "

$ws.Range("E10").Value = "  +0.06%  "

$ws.Range("D11").Value = "3.509.69"
$ws.Range("E11").Value = "  +3.94%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.08%  "

$ws.Range("E13").Value = "  -0.49%  "

$ws.Range("D14").Value = "98.185.61"
$ws.Range("E14").Value = "  +0.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.92%  "

$ws.Range("D16").Value = "4.172.57"
$ws.Range("E16").Value = "  +4.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000262"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.83%  "

$ws.Range("D19").Value = "3.506.58"
$ws.Range("E19").Value = "  +4.03%  "

$ws.Range("E20").Value = "  +9.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.85%  "

$ws.Range("E22").Value = "  -8.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "524.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.17%  "

$ws.Range("E25").Value = "  +1.09%  "

$ws.Range("E26").Value = "  +5.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "98.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.92%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.50%  "

$ws.Range("D29").Value = "3.697.61"
$ws.Range("E29").Value = "  +3.66%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +11.59%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.18%  "

$ws.Range("E32").Value = "  -2.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.192"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.597"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "30.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.994"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.02%  "

$ws.Range("E38").Value = "  +2.79%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.56%  "

$ws.Range("E40").Value = "  +3.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "532.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.69%  "

$ws.Range("E42").Value = "  +0.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.911"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.43%  "

$ws.Range("E44").Value = "  +4.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0436"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.99%  "

$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.32%  "

$ws.Range("B49").Value = "MantraDAO"
$ws.Range("C49").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.13%  "

$ws.Range("E50").Value = "  +7.93%  "

$ws.Range("E51").Value = "  +1.82%  "
